# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets to reflect the latest
# handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-25 12:17:50"
$wsZhCn.Range("H2").Value = "2016-03-25 12:18:18"
$wsZhCn.Range("E5").Value = "2016-03-25 12:17:50"
$wsZhCn.Range("H5").Value = "2016-03-25 12:18:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-25 12:17:55"
$wsDeDe.Range("H2").Value = "2016-03-25 12:18:26"
$wsDeDe.Range("E5").Value = "2016-03-25 12:17:55"
$wsDeDe.Range("H5").Value = "2016-03-25 12:18:26"
